# Apply updated crypto price/volume data per the GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.796.63"
$ws.Range("E2").Value = "  +4.09%  "
$ws.Range("D3").Value = "3.252.24"
$ws.Range("E3").Value = "  +2.20%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'396.07"
$ws.Range("E5").Value = "  -1.30%  "
$ws.Range("D6").Value = "'108.90"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "'0.580"
$ws.Range("E7").Value = "  +5.43%  "
$ws.Range("D8").Value = "3.248.81"
$ws.Range("E8").Value = "  +2.16%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").Value = "'0.625"
$ws.Range("E10").Value = "  +1.02%  "
$ws.Range("D11").Value = "'39.25"
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").Value = "'0.0972"
$ws.Range("E12").Value = "  +10.34%  "
$ws.Range("D13").Value = "'0.143"
$ws.Range("E13").Value = "  +2.32%  "
$ws.Range("D14").Value = "3.763.90"
$ws.Range("E14").Value = "  +2.44%  "
$ws.Range("D15").Value = "'8.34"
$ws.Range("E15").Value = "  +3.71%  "
$ws.Range("D16").Value = "'19.16"
$ws.Range("E16").Value = "  +0.67%  "
$ws.Range("D17").Value = "3.251.96"
$ws.Range("E17").Value = "  +2.32%  "
$ws.Range("E18").Value = "  -3.42%  "
$ws.Range("D19").Value = "'10.67"
$ws.Range("E19").Value = "  +1.69%  "
$ws.Range("D20").Value = "56.763.49"
$ws.Range("E20").Value = "  +4.33%  "
$ws.Range("E21").Value = "  +1.20%  "
$ws.Range("E22").Value = "  +9.20%  "
$ws.Range("D23").Value = "'12.92"
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("D24").Value = "'295.40"
$ws.Range("E24").Value = "  +7.49%  "
$ws.Range("D25").Value = "'74.17"
$ws.Range("E25").Value = "  +2.88%  "
$ws.Range("D26").Value = "'3.17"
$ws.Range("E26").Value = "  -3.00%  "
$ws.Range("D27").Value = "'28.21"
$ws.Range("E27").Value = "  +1.64%  "
$ws.Range("E28").Value = "  +0.94%  "
$ws.Range("E29").Value = "  -5.23%  "
$ws.Range("D30").Value = "'7.25"
$ws.Range("E30").Value = "  -2.43%  "
$ws.Range("D31").Value = "'0.168"
$ws.Range("E31").Value = "  -1.16%  "
$ws.Range("D33").Value = "'11.30"
$ws.Range("E33").Value = "  +2.31%  "
$ws.Range("D34").Value = "'0.109"
$ws.Range("E34").Value = "  -4.05%  "
$ws.Range("D35").Value = "'39.75"
$ws.Range("E35").Value = "  +7.14%  "
$ws.Range("D36").Value = "'0.0484"
$ws.Range("E37").Value = "  +2.20%  "
$ws.Range("D38").Value = "'51.52"
$ws.Range("E38").Value = "  +1.33%  "
$ws.Range("D39").Value = "'0.999"
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("E40").Value = "  -4.52%  "
$ws.Range("D41").Value = "'2.87"
$ws.Range("E41").Value = "  +1.22%  "
$ws.Range("D42").Value = "'136.83"
$ws.Range("E42").Value = "  +4.77%  "
$ws.Range("E43").Value = "  +3.65%  "
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").Value = "'3.97"
$ws.Range("E44").Value = "  -4.95%  "
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").Value = "'1.89"
$ws.Range("E45").Value = "  -2.60%  "
$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D46").Value = "'17.01"
$ws.Range("E46").Value = "  -1.87%  "
$ws.Range("D47").Value = "'0.279"
$ws.Range("E47").Value = "  -4.02%  "
$ws.Range("D48").Value = "'22.28"
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("D49").Value = "'2.13"
$ws.Range("E49").Value = "  +3.02%  "
$ws.Range("D50").Value = "2.153.99"
$ws.Range("E50").Value = "  +3.02%  "
$ws.Range("E51").Value = "  -5.44%  "
